# Update countries & provincias Spain
# This script updates the daily COVID case-count table on sheet "Pais".
# The table (rows 4..) is kept sorted descending by column B (Casos totales).
# A handful of countries received updated figures, and three of them
# ("Serbia", "Togo", "Malaui") moved up one position past their neighbour
# ("Banglades", "Zambia", "Eritrea" respectively), which kept its own
# figures unchanged but shifted down one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-Row($Row, $Country, $B, $C, $D, $E, $F, $G, $H) {
    $ws.Cells.Item($Row, 1).Value2 = $Country
    $ws.Cells.Item($Row, 2).Value2 = $B
    $ws.Cells.Item($Row, 3).Value2 = $C
    $ws.Cells.Item($Row, 4).Value2 = $D
    $ws.Cells.Item($Row, 5).Value2 = $E
    $ws.Cells.Item($Row, 6).Value2 = $F
    $ws.Cells.Item($Row, 7).Value2 = $G
    $ws.Cells.Item($Row, 8).Value2 = $H
}

# Row 4 - Estados Unidos: refreshed totals (F unchanged)
Set-Row 4 "Estados Unidos" 1162049 1275 173910 920647 16475 48 67492

# Rows 42/43 - Serbia overtakes Banglades
Set-Row 42 "Serbia" 9464 102 1551 7720 54 4 193
Set-Row 43 "Banglades" 9455 665 177 9101 1 2 177

# Row 75 - Azerbaiyan: refreshed totals (F, G unchanged)
Set-Row 75 "Azerbaiyan" 1932 38 1441 466 17 0 25

# Row 127 - Mauricio: refreshed totals (B, C, F, G, H unchanged)
Set-Row 127 "Mauricio" 332 0 315 7 3 0 10

# Rows 147/148 - Togo overtakes Zambia
Set-Row 147 "Togo" 124 1 67 48 0 0 9
Set-Row 148 "Zambia" 124 5 78 43 1 0 3

# Rows 176/177 - Malaui overtakes Eritrea
Set-Row 176 "Malaui" 39 1 9 27 1 0 3
Set-Row 177 "Eritrea" 39 0 26 13 0 0 0
